$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 so the existing rows (years
# 2008-2024) shift down one, making room for a new "2007" leading row.
$ws.Rows.Item(2).Insert()

# The inserted row picks up the formatting of the row above it (the bold
# header row), which is wrong for a plain data row - clear that first.
$ws.Range("B2:E2").ClearFormats()
$ws.Range("E2").ClearContents()

# Column A keeps its date number-format/border style on every data row;
# copy that formatting (only) from the row below onto the new A2 cell.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Full target table (rows 2-19): A=serial date, B=y0 year, C=y0_forecast,
# D=y1 year, E=y1_forecast (E intentionally blank for the first three rows).
$data = @(
  @(2, 39400, 2007, 5.896808312953783, 2008, $null),
  @(3, 39765, 2008, 7.441962824572235, 2009, $null),
  @(4, 40130, 2009, 6.277541464866987, 2010, $null),
  @(5, 40494, 2010, 6.535114773304773, 2011, 6.325696408067327),
  @(6, 40862, 2011, 5.12051970717502, 2012, 4.950888348161886),
  @(7, 41228, 2012, 3.65682115264816, 2013, 3.982564147794321),
  @(8, 41592, 2013, 2.943878639034381, 2014, 4.334309403335435),
  @(9, 41957, 2014, 1.172679597477866, 2015, 2.644356903452572),
  @(10, 42321, 2015, 2.961845079861303, 2016, 3.383932287548697),
  @(11, 42689, 2016, 2.508469427909898, 2017, 3.355044026998955),
  @(12, 43053, 2017, 3.523703831572056, 2018, 3.74984170812418),
  @(13, 43418, 2018, 1.178605266817589, 2019, 2.186196327763934),
  @(14, 43783, 2019, 3.047037961814492, 2020, 2.880436144359444),
  @(15, 44159, 2020, -0.2228847697281378, 2021, 1.982741503124119),
  @(16, 44525, 2021, -1.165854108406617, 2022, 2.782217648649521),
  @(17, 44890, 2022, 2.501311189006916, 2023, 2.985901060752827),
  @(18, 45254, 2023, 0.6753076481029074, 2024, 0.7957830962485257),
  @(19, 45618, 2024, 2.039329803030121, 2025, 2.510359031091491)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  if ($row[5] -ne $null) {
    $ws.Cells.Item($r, 5).Value = $row[5]
  }
}
